# Re-commit the "D suite.xlsx" test-results corrections on the "Test Cases"
# sheet (the previously committed copy had gotten corrupted).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Column E ("Results"): rows 2 and 3 were blank, now recorded as SKIP.
$ws.Range("E2").Value = "SKIP"
$ws.Range("E3").Value = "SKIP"

# Rows 6, 7 and 9 were (incorrectly) recorded as FAIL; correct them to SKIP.
$ws.Range("E6").Value = "SKIP"
$ws.Range("E7").Value = "SKIP"
$ws.Range("E9").Value = "SKIP"

# Row 8 was (incorrectly) recorded as PASS; correct it to SKIP.
$ws.Range("E8").Value = "SKIP"

# Row 41 was recorded as SKIP; correct it to PASS.
$ws.Range("E41").Value = "PASS"

# Refresh the sheet view/selection state that Excel persisted for this sheet.
$ws.Activate()
$ws.Range("C34").Select()
